$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the two measurement comment cells with revised values ---
# G7: delta values updated (-0.17MW -> -0.15MW)
$ws.Range("G7").Value2 = "Δt = 1.69s, ΔPmoy = -0.15MW, ΔPmax = -0.04MW "
# F7: recharge test values updated
$ws.Range("F7").Value2 = "t=0.31s, Pmoy = 2.55MW et Pmax=3.62MW pour PSIM, Pmoy=2.48MW et Pmax=3.56MW pour SPS"

# --- Row height adjustments (rows grew taller, likely to fit inserted images) ---
$ws.Rows.Item(2).RowHeight = 51
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 109.5
$ws.Rows.Item(8).RowHeight = 102.75
$ws.Rows.Item(9).RowHeight = 37.5
$ws.Rows.Item(10).RowHeight = 58.5
$ws.Rows.Item(12).RowHeight = 37.5

# --- Update the view selection to match the final state ---
$ws.Activate()
$ws.Range("K9").Select()
